$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column I for "Left Wetted" / "Right Wetted"
$ws.Columns("I:J").Insert()

# New column headers
$ws.Range("I1").Value = "Left Wetted"
$ws.Range("J1").Value = "Right Wetted"

# Fix transcription error: B5 was a flat 6.5, should be 6.5 feet converted to inches (6.5 * 12)
$ws.Range("B5").Formula = "= 6.5 * 12"

# Fix transcription error: H10 was 45, should be 4.5
$ws.Range("H10").Value = 4.5

# Restore the selection/view state left behind by the authoring tool
$ws.Range("Q15").Select()
